$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range('D2').Value = '62.324.61'
$ws.Range('E2').Value = '  -1.99%  '
$ws.Range('D3').Value = '2.432.76'
$ws.Range('E3').Value = '  -1.63%  '
$ws.Range('D5').NumberFormat = '@'
$ws.Range('D5').Value = '570.25'
$ws.Range('D5').Style = 'Normal'
$ws.Range('E5').Value = '  -1.19%  '
$ws.Range('D6').NumberFormat = '@'
$ws.Range('D6').Value = '143.47'
$ws.Range('D6').Style = 'Normal'
$ws.Range('E6').Value = '  -3.60%  '
$ws.Range('E7').Value = '  +0.16%  '
$ws.Range('E8').Value = '  -2.29%  '
$ws.Range('D9').Value = '2.428.91'
$ws.Range('E9').Value = '  -2.09%  '
$ws.Range('E10').Value = '  -5.68%  '
$ws.Range('E11').Value = '  +0.53%  '
$ws.Range('E12').Value = '  -2.26%  '
$ws.Range('E13').Value = '  -3.21%  '
$ws.Range('D14').NumberFormat = '@'
$ws.Range('D14').Value = '26.59'
$ws.Range('D14').Style = 'Normal'
$ws.Range('E14').Value = '  -2.36%  '
$ws.Range('D15').NumberFormat = '@'
$ws.Range('D15').Value = '0.0000173'
$ws.Range('D15').Style = 'Normal'
$ws.Range('E15').Value = '  -6.23%  '
$ws.Range('D16').Value = '2.876.93'
$ws.Range('E16').Value = '  -1.47%  '
$ws.Range('D17').Value = '62.238.87'
$ws.Range('E17').Value = '  -1.91%  '
$ws.Range('D18').Value = '2.429.89'
$ws.Range('E18').Value = '  -1.97%  '
$ws.Range('D19').NumberFormat = '@'
$ws.Range('D19').Value = '11.04'
$ws.Range('D19').Style = 'Normal'
$ws.Range('E19').Value = '  -4.29%  '
$ws.Range('E20').Value = '  -3.20%  '
$ws.Range('D21').NumberFormat = '@'
$ws.Range('D21').Value = '324.54'
$ws.Range('D21').Style = 'Normal'
$ws.Range('E21').Value = '  -1.38%  '
$ws.Range('E22').Value = '  -2.75%  '
$ws.Range('E23').Value = '  +2.22%  '
$ws.Range('E24').Value = '  +0.27%  '
$ws.Range('D25').NumberFormat = '@'
$ws.Range('D25').Value = '65.07'
$ws.Range('D25').Style = 'Normal'
$ws.Range('D26').NumberFormat = '@'
$ws.Range('D26').Value = '620.10'
$ws.Range('D26').Style = 'Normal'
$ws.Range('E26').Value = '  -2.42%  '
$ws.Range('D27').NumberFormat = '@'
$ws.Range('D27').Value = '9.06'
$ws.Range('D27').Style = 'Normal'
$ws.Range('E27').Value = '  +1.20%  '
$ws.Range('E28').Value = '  -9.59%  '
$ws.Range('E29').Value = '  -1.67%  '
$ws.Range('E30').Value = '  +0.37%  '
$ws.Range('E31').Value = '  -5.56%  '
$ws.Range('D32').NumberFormat = '@'
$ws.Range('D32').Value = '8.02'
$ws.Range('D32').Style = 'Normal'
$ws.Range('E32').Value = '  -4.69%  '
$ws.Range('E33').Value = '  -3.31%  '
$ws.Range('E34').Value = '  -8.00%  '
$ws.Range('D35').NumberFormat = '@'
$ws.Range('D35').Value = '5.03'
$ws.Range('D35').Style = 'Normal'
$ws.Range('E35').Value = '  -3.46%  '
$ws.Range('E36').Value = '  +0.32%  '
$ws.Range('E37').Value = '  -6.94%  '
$ws.Range('E38').Value = '  -3.12%  '
$ws.Range('D39').NumberFormat = '@'
$ws.Range('D39').Value = '18.55'
$ws.Range('D39').Style = 'Normal'
$ws.Range('E39').Value = '  -2.42%  '
$ws.Range('D40').NumberFormat = '@'
$ws.Range('D40').Value = '147.17'
$ws.Range('D40').Style = 'Normal'
$ws.Range('E40').Value = '  +0.32%  '
$ws.Range('E41').Value = '  -5.83%  '
$ws.Range('E42').Value = '  -7.14%  '
$ws.Range('D43').NumberFormat = '@'
$ws.Range('D43').Value = '42.52'
$ws.Range('D43').Style = 'Normal'
$ws.Range('E43').Value = '  +1.52%  '
$ws.Range('E45').Value = '  -8.43%  '
$ws.Range('D46').NumberFormat = '@'
$ws.Range('D46').Value = '145.17'
$ws.Range('D46').Style = 'Normal'
$ws.Range('E46').Value = '  -3.74%  '
$ws.Range('D47').NumberFormat = '@'
$ws.Range('D47').Value = '3.69'
$ws.Range('D47').Style = 'Normal'
$ws.Range('E47').Value = '  -2.39%  '
$ws.Range('B48').Value = 'InjectiveProtocol'
$ws.Range('C48').Value = 'https://coinranking.com/coin/PkY9BmsyW+injectiveprotocol-inj'
$ws.Range('D48').NumberFormat = '@'
$ws.Range('D48').Value = '20.11'
$ws.Range('D48').Style = 'Normal'
$ws.Range('E48').Value = '  -4.93%  '
$ws.Range('B49').Value = 'Hedera'
$ws.Range('C49').Value = 'https://coinranking.com/coin/jad286TjB+hedera-hbar'
$ws.Range('D49').NumberFormat = '@'
$ws.Range('D49').Value = '0.0521'
$ws.Range('D49').Style = 'Normal'
$ws.Range('E49').Value = '  -5.41%  '
$ws.Range('D50').NumberFormat = '@'
$ws.Range('D50').Value = '0.593'
$ws.Range('D50').Style = 'Normal'
$ws.Range('E50').Value = '  -2.92%  '
$ws.Range('E51').Value = '  -4.65%  '
